# Edit: rename header "value" -> "first_release_value", and replace the A:B
# time-series data with an extended/refreshed dataset (53 -> 84 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header text change (B1 shared string "value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# 2) Extend the data block with new rows 54-84, carrying over the same
#    formatting (date number format etc.) that rows 2-53 already use.
$ws.Range("A2:B2").Copy()
$ws.Range("A54:B84").PasteSpecial(-4122)

# 3) Write the refreshed date/value series into A2:B84.
$dates = @(38398, 38487, 38579, 38671, 38763, 38852, 38944, 39036, 39128, 39217, 39309, 39401, 39493, 39583, 39675, 39767, 39859, 39948, 40040, 40132, 40224, 40313, 40405, 40497, 40589, 40678, 40770, 40862, 40954, 41044, 41136, 41228, 41320, 41409, 41501, 41593, 41685, 41774, 41866, 41958, 42050, 42139, 42231, 42323, 42415, 42505, 42597, 42689, 42781, 42870, 42962, 43054, 43146, 43235, 43327, 43419, 43511, 43600, 43692, 43784, 43876, 43966, 44058, 44150, 44242, 44331, 44423, 44515, 44607, 44696, 44788, 44880, 44972, 45061, 45153, 45245, 45337, 45427, 45519, 45611, 45703, 45792, 45884)
$vals = @(-1.480585662758955, 4.695917448367524, 3.170055072151484, 2.606847295028047, 3.777024304539765, 0.7137374907328251, 3.543396901532475, 1.324703771913761, 2.324224946277042, -1.781652582373326, 3.94593981855364, 0.2255418731353984, 3.488038255381227, -1.392382908151674, 3.815959839717081, -4.107770248796484, -5.35640370103539, -5.410562843974105, 4.973955294124409, -1.616494377065351, 6.089842363966454, 7.771920357185309, 1.933642100495049, 0.796738168115894, 1.463284815360069, 2.891754278273439, 2.64337718803263, -0.7880943707230443, 0.02227176351210858, 2.238969365349575, 1.043875137114455, -1.312661112537128, -2.1, 1.933078912701916, 0.8117228711496978, 1.285328021680314, 2.2085072997628, 1.119204613350774, 1.677790477509291, 1.857496130824472, 1.464045137806849, 0.4869808267284412, 1.144719845809078, 0.5201842158159025, 1.396417402226163, 0.09677264461834056, 0.1611306858251567, 2.489390679284554, 0.3983205376114825, 2.396748302637434, 0.8920379051669016, 1.989289785701104, -1.095080621818852, 1.709973952921786, 1.3, 0.7345957212796748, 0.7, -0.2740143521242828, 0.06866907528606703, 1.271957369209446, -1.6, -16.02569689670956, 7.96955251685678, 2.117788110998191, -0.9861240056009706, 2.876944405321424, 0.4364757668776207, 0.6542354095451515, -1.358640149334988, 1.767346889326234, 0.1973819540654631, -2.016521230865749, -1.429923541452922, 0.1801827438520291, 0.6692453970872521, -1.33243152085096, -0.09168777270478756, 0.8481739611978583, -0.2604183589432552, 0.3093370292089048, -0.03180050048325711, 0.04735640278761366, 0.2461857363876589)

for ($i = 0; $i -lt $dates.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $dates[$i]
  $ws.Cells.Item($row, 2).Value = $vals[$i]
}
